$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ARCHITECTURE")

# Update ventilation property (column H, n50 -> win_op) from 0.9 to 0.5 for rows 26-181
for ($r = 26; $r -le 181; $r++) {
    $ws.Cells.Item($r, 8).Value = 0.5
}

# Activate the sheet and set the selection to match the saved view state
$ws.Activate()
$ws.Range("H26:H181").Select()
